$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.409.72"
$ws.Range("E2").Value = "  +2.08%  "

$ws.Range("D3").Value = "3.171.60"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.95%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.169.94"
$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("E10").Value = "  +1.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.75%  "

$ws.Range("D15").Value = "3.692.83"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("D17").Value = "3.171.46"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").Value = "63.363.71"
$ws.Range("E18").Value = "  +1.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.697"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.92%  "

$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("E29").Value = "  +4.09%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.94%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("E33").Value = "  -0.97%  "

$ws.Range("E34").Value = "  +2.68%  "

$ws.Range("E35").Value = "  -1.27%  "

$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("D37").Value = "0.0₃0733"
$ws.Range("E37").Value = "  +6.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0390"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.84%  "

$ws.Range("E40").Value = "  +1.76%  "

$ws.Range("E41").Value = "  +0.83%  "

$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "392.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.45%  "

$ws.Range("D44").Value = "2.806.17"
$ws.Range("E44").Value = "  -4.66%  "

$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.64%  "

$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.23"
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = "  +0.99%  "
